$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C, shifting existing "Number Of Options" .. "Answer" columns
# one place to the right, and populate the new header cell.
$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "Question Type"

# Match the new column widths: B & C share one width, D gets the other.
$ws.Range("B1:C1").ColumnWidth = 31.498697916666668
$ws.Range("D1").ColumnWidth = 17.498697916666668

# Update the active selection to reflect where the editor left off.
$ws.Range("E23").Select()
